# Insert a new weekly price record as row 202, pushing existing rows 202-221
# down to 203-222 (the sheet keeps a daily/weekly log of price observations,
# newest entries are inserted near the top of the historical block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 202; rows below shift down by one.
$ws.Rows("202:202").Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(202, 1).Value = 2
$ws.Cells.Item(202, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(202, 3).Value = "Coquimbo"
$ws.Cells.Item(202, 4).Value = 45106
$ws.Cells.Item(202, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(202, 5).Value = 4
$ws.Cells.Item(202, 6).Value = 100112043
$ws.Cells.Item(202, 7).Value = "Pepino ensalada"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 500
$ws.Cells.Item(202, 11).Value = 11000
$ws.Cells.Item(202, 12).Value = 12000
$ws.Cells.Item(202, 13).Value = 11500
$ws.Cells.Item(202, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(202, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(202, 16).Value = 164
$ws.Cells.Item(202, 17).Value = 70
$ws.Cells.Item(202, 18).Value = "Hortaliza"
